# Trade #24 closed at 2026-02-16 22:54:28 - base_strategy UP +0.000%
# Append the new trade row (row 25) to both the "All Trades" and
# "base_strategy" worksheets - they track the same trade log.

$wb = $excel.ActiveWorkbook

$sheetNames = @("All Trades", "base_strategy")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    $row = 25

    $ws.Range("A$row").Value = 24

    # Force the date-looking string to stay plain text (matches the rest
    # of column B) instead of being auto-converted into a date serial.
    $ws.Range("B$row").NumberFormat = "@"
    $ws.Range("B$row").Value = "2026-02-16"
    $ws.Range("B$row").Style = "Normal"

    $ws.Range("C$row").Value = "22:54:28"
    $ws.Range("D$row").Value = "base_strategy"
    $ws.Range("E$row").Value = "UP"
    $ws.Range("F$row").Value = 49.999998
    $ws.Range("H$row").Value = "OPEN"
    $ws.Range("I$row").Value = 0
    $ws.Range("J$row").Value = 0
    $ws.Range("K$row").Value = 100
    $ws.Range("L$row").Value = 0
    $ws.Range("M$row").Value = 0
    $ws.Range("N$row").Value = 0.6
    $ws.Range("O$row").Value = "Normal spread capture: 19600 bps"
    $ws.Range("Q$row").Value = 0
}
